# Update Lrpap1-Sort1 NATMI TPM values (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 4.573811
$ws.Range("H2").Value = 13.721433
$ws.Range("I2").Value = 0.1659009079913533
$ws.Range("J2").Value = 0.1659009079913533
$ws.Range("M2").Value = 0.4652636666666667
$ws.Range("N2").Value = 1.395791
$ws.Range("O2").Value = 0.02604271297411062
$ws.Range("P2").Value = 0.02604271297411062
$ws.Range("Q2").Value = 2.128028076500334
$ws.Range("R2").Value = 19.152252688503
$ws.Range("S2").Value = 0.00432050972896315
$ws.Range("T2").Value = 0.004320509728963149

# Row 3
$ws.Range("G3").Value = 4.573811
$ws.Range("H3").Value = 13.721433
$ws.Range("I3").Value = 0.1659009079913533
$ws.Range("J3").Value = 0.1659009079913533
$ws.Range("O3").Value = 0.09971126509087273
$ws.Range("P3").Value = 0.09971126509087272
$ws.Range("Q3").Value = 8.147706111405666
$ws.Range("R3").Value = 73.329355002651
$ws.Range("S3").Value = 0.01654218941554232
$ws.Range("T3").Value = 0.01654218941554231

# Row 4
$ws.Range("G4").Value = 4.573811
$ws.Range("H4").Value = 13.721433
$ws.Range("I4").Value = 0.1659009079913533
$ws.Range("J4").Value = 0.1659009079913533
$ws.Range("M4").Value = 15.618761
$ws.Range("N4").Value = 46.856283
$ws.Range("O4").Value = 0.8742460219350168
$ws.Range("P4").Value = 0.8742460219350167
$ws.Range("Q4").Value = 71.437260868171
$ws.Range("R4").Value = 642.935347813539
$ws.Range("S4").Value = 0.1450382088468479
$ws.Range("T4").Value = 0.1450382088468478

# Row 5
$ws.Range("I5").Value = 0.5322852674812913
$ws.Range("J5").Value = 0.5322852674812913
$ws.Range("M5").Value = 0.4652636666666667
$ws.Range("N5").Value = 1.395791
$ws.Range("O5").Value = 0.02604271297411062
$ws.Range("P5").Value = 0.02604271297411062
$ws.Range("Q5").Value = 6.827678085804779
$ws.Range("R5").Value = 61.449102772243
$ws.Range("S5").Value = 0.01386215244136297
$ws.Range("T5").Value = 0.01386215244136297

# Row 6
$ws.Range("I6").Value = 0.5322852674812913
$ws.Range("J6").Value = 0.5322852674812913
$ws.Range("O6").Value = 0.09971126509087273
$ws.Range("P6").Value = 0.09971126509087272
$ws.Range("Q6").Value = 26.14153219158122
$ws.Range("R6").Value = 235.273789724231
$ws.Range("S6").Value = 0.05307483740979314
$ws.Range("T6").Value = 0.05307483740979314

# Row 7
$ws.Range("I7").Value = 0.5322852674812913
$ws.Range("J7").Value = 0.5322852674812913
$ws.Range("M7").Value = 15.618761
$ws.Range("N7").Value = 46.856283
$ws.Range("O7").Value = 0.8742460219350168
$ws.Range("P7").Value = 0.8742460219350167
$ws.Range("Q7").Value = 229.2030946046843
$ws.Range("R7").Value = 2062.827851442159
$ws.Range("S7").Value = 0.4653482776301353
$ws.Range("T7").Value = 0.4653482776301353

# Row 8
$ws.Range("G8").Value = 8.320867
$ws.Range("H8").Value = 24.962601
$ws.Range("I8").Value = 0.3018138245273554
$ws.Range("J8").Value = 0.3018138245273554
$ws.Range("M8").Value = 0.4652636666666667
$ws.Range("N8").Value = 1.395791
$ws.Range("O8").Value = 0.02604271297411062
$ws.Range("P8").Value = 0.02604271297411062
$ws.Range("Q8").Value = 3.871397090265667
$ws.Range("R8").Value = 34.842573812391
$ws.Range("S8").Value = 0.007860050803784506
$ws.Range("T8").Value = 0.007860050803784503

# Row 9
$ws.Range("G9").Value = 8.320867
$ws.Range("H9").Value = 24.962601
$ws.Range("I9").Value = 0.3018138245273554
$ws.Range("J9").Value = 0.3018138245273554
$ws.Range("O9").Value = 0.09971126509087273
$ws.Range("P9").Value = 0.09971126509087272
$ws.Range("Q9").Value = 14.82264547181633
$ws.Range("R9").Value = 133.403809246347
$ws.Range("S9").Value = 0.03009423826553728
$ws.Range("T9").Value = 0.03009423826553728

# Row 10
$ws.Range("G10").Value = 8.320867
$ws.Range("H10").Value = 24.962601
$ws.Range("I10").Value = 0.3018138245273554
$ws.Range("J10").Value = 0.3018138245273554
$ws.Range("M10").Value = 15.618761
$ws.Range("N10").Value = 46.856283
$ws.Range("O10").Value = 0.8742460219350168
$ws.Range("P10").Value = 0.8742460219350167
$ws.Range("Q10").Value = 129.961632985787
$ws.Range("R10").Value = 1169.654696872083
$ws.Range("S10").Value = 0.2638595354580336
$ws.Range("T10").Value = 0.2638595354580336
